$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-8: columns D, M, N, O, P, S
$ws.Range("D2").Value = 44216
$ws.Range("M2").Value = 55
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 11545
$ws.Range("S2").Value = 825

$ws.Range("D3").Value = 44253
$ws.Range("M3").Value = 90
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 12667
$ws.Range("S3").Value = 905

$ws.Range("D4").Value = 44181
$ws.Range("M4").Value = 65
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 9462
$ws.Range("S4").Value = 676

$ws.Range("D5").Value = 44172
$ws.Range("M5").Value = 90
$ws.Range("N5").Value = 8500
$ws.Range("O5").Value = 9000
$ws.Range("P5").Value = 8806
$ws.Range("S5").Value = 629

$ws.Range("D6").Value = 44210
$ws.Range("M6").Value = 70
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 11000
$ws.Range("P6").Value = 10357
$ws.Range("S6").Value = 740

$ws.Range("D7").Value = 44232
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 11583
$ws.Range("S7").Value = 827

$ws.Range("D8").Value = 44229
$ws.Range("M8").Value = 55
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 11364
$ws.Range("S8").Value = 812

# Add new row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 45138
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100102
$ws.Range("H9").Value = "Cítricos"
$ws.Range("I9").Value = 100102006
$ws.Range("J9").Value = "Pomelo"
$ws.Range("K9").Value = "Start Ruby"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 14000
$ws.Range("P9").Value = 14000
$ws.Range("Q9").Value = "$/caja 14 kilos granel"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1000
$ws.Range("T9").Value = 14

# Apply the same style as D2:D8 (date format) to D9
$ws.Range("D2").Copy()
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Value = 45138
